# "Conclusão da implementação do gamma"
# Item "1s" (row 26) is now finished: progress goes from 25% to 100%,
# which in turn flips its Status formula result to "Done!" and moves the
# overall average in I3 upward (both cells are formula-driven and
# recalculate automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D26").Value = 100

# Leave the selection where the author ended up after the edit.
[void]$ws.Range("D27").Select()
